$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the lookup table (header + 3 rows) ------------------------
# Order matters: it controls the order new entries land in sharedStrings.xml.
# Target shared-string order is: principle, teacher, hr, role, id
$ws.Range("B2").Value = "principle"
$ws.Range("B3").Value = "teacher"
$ws.Range("B4").Value = "hr"
$ws.Range("B1").Value = "role"
$ws.Range("A1").Value = "id"

# Ids now start at 1 instead of 0
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- New column C: generated SQL INSERT statements ----------------------
$ws.Range("C2").Formula = "=CONCAT(""INSERT INTO permission (role) VALUE ('"",B2,""');"")"
$ws.Range("C3").Formula = "=CONCAT(""INSERT INTO permission (role) VALUE ('"",B3,""');"")"
$ws.Range("C4").Formula = "=CONCAT(""INSERT INTO permission (role) VALUE ('"",B4,""');"")"

# --- Formatting ----------------------------------------------------------
# Left-align every used cell (A1:C4), including the otherwise-empty C1.
$ws.Range("A1:C4").HorizontalAlignment = -4131  ## xlLeft

# Auto-size column C to fit the generated SQL text.
$ws.Columns.Item(3).AutoFit() | Out-Null

# --- Selection / view state ----------------------------------------------
$ws.Range("J9").Select() | Out-Null
